$wb = $excel.ActiveWorkbook

# Rename "prok only" sheet to "DB peps" (sheet2)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "DB peps"

# Update the active selection on sheet2: drop topLeftCell scroll anchor,
# move the active cell/selection to D11
[void]$ws2.Range("D11").Select()

# Increase the height of row 3 on sheet2 from 33 to 35
$ws2.Rows.Item(3).RowHeight = 35
